$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1.721470851391172, 0.9234690187349804),
    @(3, 3.486192098951822, 1.944776428788816),
    @(4, 5.291008189270382, 3.169215447447818),
    @(5, 5.587584589080254, 4.090188008681169),
    @(6, 10.79464947974671, 5.042048921878017),
    @(7, 11.61253535479567, 5.900352457760341),
    @(8, 12.02696128401184, 7.395425997841303),
    @(9, 15.12996035942051, 8.347768715480118),
    @(10, 16.53251100512178, 9.346900283863405),
    @(11, 16.81680184152511, 10.19490799257818),
    @(12, 17.43491540630607, 11.00354207578137),
    @(13, 18.31605758012822, 12.0660174305829),
    @(14, 21.29847252930778, 13.39321972155926),
    @(15, 26.06078701724958, 14.61595394567654),
    @(16, 27.42617037529666, 15.62119540484504),
    @(17, 27.6750834583812, 16.54262382415783),
    @(18, 27.84690271063112, 17.57604501832896),
    @(19, 30.74247008811135, 19.59722326389782),
    @(20, 31.38430989350973, 20.60096842483541),
    @(21, 33.4119129860735, 21.67735797213086),
    @(22, 37.36764397841052, 22.67995873243806),
    @(23, 39.44914260255963, 23.59435296827072),
    @(24, 39.77514401976434, 24.9234581430823),
    @(25, 39.93362051959899, 25.81863002316024),
    @(26, 43.67401244288455, 26.80947622745665),
    @(27, 47.81403792188777, 27.72333735057278),
    @(28, 47.97811568950357, 28.7825416204646),
    @(29, 52.81887057730511, 30.31721839728307),
    @(30, 53.07783868975601, 31.49556251505136),
    @(31, 53.90148682097326, 32.41442166081525),
    @(32, 58.00674136434694, 33.3409564680513),
    @(33, 64.87479573449659, 34.27945796942749),
    @(34, 65.33836825404471, 35.45960740508962),
    @(35, 68.2367885623102, 36.64841443579079),
    @(36, 70.5891350281524, 37.59966142548254),
    @(37, 72.8914300464168, 38.77413408081087),
    @(38, 73.20859702792586, 39.91361350474125),
    @(39, 77.26228268374241, 41.12135812766236),
    @(40, 78.90674977238301, 42.02951212724357),
    @(41, 79.96219785761541, 43.2011265791078),
    @(42, 82.26516207262878, 44.5205962535062),
    @(43, 82.39474874932522, 45.34907153737197),
    @(44, 83.07104152212452, 46.52626588052185),
    @(45, 83.57871486788777, 47.71808363082739),
    @(46, 83.78785685756327, 49.04603829320707),
    @(47, 84.04936068813316, 49.89280821882973),
    @(48, 87.70181312624476, 51.07171699501184),
    @(49, 89.84885690733311, 52.35562049597146),
    @(50, 96.78473914802133, 53.55599024731602)
)

foreach ($row in $data) {
    $r = $row[0]
    $bVal = $row[1]
    $cVal = $row[2]
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
}
